$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 273.3889
$ws.Range("I19").Value = 327.625
$ws.Range("J19").Value = 230
$ws.Range("K19").Value = 327.625
$ws.Range("L19").Value = 230
$ws.Range("M19").Value = -152.625
$ws.Range("N19").Value = -580

$ws.Range("H28").Value = 4089.1667
$ws.Range("I28").Value = 515.3570999999999
$ws.Range("J28").Value = 16597.5
$ws.Range("K28").Value = 515.3570999999999
$ws.Range("L28").Value = 16597.5
$ws.Range("M28").Value = -30.35709999999995
$ws.Range("N28").Value = -17567.5

$ws.Range("H53").Value = 445.54544
$ws.Range("I53").Value = 300.66666
$ws.Range("J53").Value = 499.875
$ws.Range("K53").Value = 300.66666
$ws.Range("L53").Value = 499.875
$ws.Range("M53").Value = 336.33334
$ws.Range("N53").Value = -1773.875

$ws.Range("H74").Value = 3715.32
$ws.Range("I74").Value = 3639.0625
$ws.Range("K74").Value = 3639.0625
$ws.Range("M74").Value = -2703.0625

$ws.Range("H77").Value = 3715.32
$ws.Range("I77").Value = 3639.0625
$ws.Range("K77").Value = 18195.3125
$ws.Range("M77").Value = -13515.3125

$ws.Range("H129").Value = 979.7738000000001
$ws.Range("I129").Value = 274.4
$ws.Range("J129").Value = 1024.4177
$ws.Range("K129").Value = 823.1999999999999
$ws.Range("L129").Value = 3073.2531
$ws.Range("M129").Value = 4176.8
$ws.Range("N129").Value = -13073.2531

$ws.Range("H132").Value = 2609784.2
$ws.Range("I132").Value = 3054340.2
$ws.Range("J132").Value = 5957.143
$ws.Range("K132").Value = 9163020.600000001
$ws.Range("L132").Value = 17871.429
$ws.Range("M132").Value = -9160490.600000001
$ws.Range("N132").Value = -22931.429

$ws.Range("H137").Value = 2701.1372
$ws.Range("I137").Value = 2501.4055
$ws.Range("J137").Value = 3229
$ws.Range("K137").Value = 7504.2165
$ws.Range("L137").Value = 9687
$ws.Range("M137").Value = -4954.2165
$ws.Range("N137").Value = -14787

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 102
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H61").Value = 5296.3335
$ws.Range("I61").Value = 6695
$ws.Range("J61").Value = 1100.3334
$ws.Range("K61").Value = 6695
$ws.Range("L61").Value = 1100.3334
$ws.Range("M61").Value = -6483
$ws.Range("N61").Value = -1524.3334

$ws.Range("H74").Value = 2633.08
$ws.Range("I74").Value = 3265.6487
$ws.Range("J74").Value = 832.6923
$ws.Range("K74").Value = 3265.6487
$ws.Range("L74").Value = 832.6923
$ws.Range("M74").Value = -2391.6487
$ws.Range("N74").Value = -2580.6923

$ws.Range("H77").Value = 2633.08
$ws.Range("I77").Value = 3265.6487
$ws.Range("J77").Value = 832.6923
$ws.Range("K77").Value = 16328.2435
$ws.Range("L77").Value = 4163.4615
$ws.Range("M77").Value = -11960.2435
$ws.Range("N77").Value = -12899.4615

$ws.Range("H122").Value = 1455.1904
$ws.Range("I122").Value = 1369.9445
$ws.Range("J122").Value = 1966.6666
$ws.Range("K122").Value = 4109.833500000001
$ws.Range("L122").Value = 5899.9998
$ws.Range("M122").Value = -1659.833500000001
$ws.Range("N122").Value = -10799.9998

$ws.Range("H132").Value = 6194.614
$ws.Range("I132").Value = 4849
$ws.Range("J132").Value = 9078.071
$ws.Range("K132").Value = 14547
$ws.Range("L132").Value = 27234.213
$ws.Range("M132").Value = -12017
$ws.Range("N132").Value = -32294.213

$ws.Range("H136").Value = 5296.3335
$ws.Range("I136").Value = 6695
$ws.Range("J136").Value = 1100.3334
$ws.Range("K136").Value = 20085
$ws.Range("L136").Value = 3301.0002
$ws.Range("M136").Value = -17535
$ws.Range("N136").Value = -8401.0002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 102
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H86").Value = 2976.561
$ws.Range("I86").Value = 3072.3225
$ws.Range("J86").Value = 2679.7
$ws.Range("K86").Value = 3072.3225
$ws.Range("L86").Value = 2679.7
$ws.Range("M86").Value = -1949.3225
$ws.Range("N86").Value = -4925.7

$ws.Range("H89").Value = 2976.561
$ws.Range("I89").Value = 3072.3225
$ws.Range("J89").Value = 2679.7
$ws.Range("K89").Value = 15361.6125
$ws.Range("L89").Value = 13398.5
$ws.Range("M89").Value = -9745.612500000001
$ws.Range("N89").Value = -24630.5

$ws.Range("H94").Value = 1704.9166
$ws.Range("I94").Value = 995.44446
$ws.Range("J94").Value = 3833.3333
$ws.Range("K94").Value = 995.44446
$ws.Range("L94").Value = 3833.3333
$ws.Range("M94").Value = -544.44446
$ws.Range("N94").Value = -4735.3333

$ws.Range("H99").Value = 1390.0476
$ws.Range("I99").Value = 1163.3334
$ws.Range("J99").Value = 1956.8334
$ws.Range("K99").Value = 1163.3334
$ws.Range("L99").Value = 1956.8334
$ws.Range("M99").Value = 334.6666
$ws.Range("N99").Value = -4952.8334

$ws.Range("H107").Value = 1499.125
$ws.Range("I107").Value = 931.4545000000001
$ws.Range("J107").Value = 2748
$ws.Range("K107").Value = 931.4545000000001
$ws.Range("L107").Value = 2748
$ws.Range("M107").Value = 988.5454999999999
$ws.Range("N107").Value = -6588

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1274.7
$ws.Range("I58").Value = 1447.8148
$ws.Range("J58").Value = 915.1539
$ws.Range("K58").Value = 1447.8148
$ws.Range("L58").Value = 915.1539
$ws.Range("M58").Value = -1244.8148
$ws.Range("N58").Value = -1321.1539

$ws.Range("H132").Value = 10701.954
$ws.Range("I132").Value = 8660.385
$ws.Range("J132").Value = 13650.889
$ws.Range("K132").Value = 25981.155
$ws.Range("L132").Value = 40952.667
$ws.Range("M132").Value = -23451.155
$ws.Range("N132").Value = -46012.667

$ws.Range("H136").Value = 1274.7
$ws.Range("I136").Value = 1447.8148
$ws.Range("J136").Value = 915.1539
$ws.Range("K136").Value = 4343.4444
$ws.Range("L136").Value = 2745.4617
$ws.Range("M136").Value = -1793.4444
$ws.Range("N136").Value = -7845.4617

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 325.75
$ws.Range("I8").Value = 325.75
$ws.Range("K8").Value = 977.25
$ws.Range("M8").Value = -838.25

$ws.Range("H68").Value = 15231.143
$ws.Range("J68").Value = 15231.143
$ws.Range("L68").Value = 45693.429
$ws.Range("N68").Value = -47315.429

$ws.Range("H71").Value = 15231.143
$ws.Range("J71").Value = 15231.143
$ws.Range("L71").Value = 137080.287
$ws.Range("N71").Value = -145192.287

$ws.Range("H86").Value = 0
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("M86").ClearContents()
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 0
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("M89").ClearContents()
$ws.Range("N89").ClearContents()

$ws.Range("H98").Value = 454817.12
$ws.Range("I98").Value = 257.93332
$ws.Range("J98").Value = 1428872.6
$ws.Range("K98").Value = 773.7999599999999
$ws.Range("L98").Value = 4286617.800000001
$ws.Range("M98").Value = 724.2000400000001
$ws.Range("N98").Value = -4289613.800000001

$ws.Range("H132").Value = 1010.5714
$ws.Range("I132").Value = 627
$ws.Range("J132").Value = 1298.25
$ws.Range("K132").Value = 5643
$ws.Range("L132").Value = 11684.25
$ws.Range("M132").Value = -3113
$ws.Range("N132").Value = -16744.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1237.12
$ws.Range("I122").Value = 1031.3684
$ws.Range("J122").Value = 1888.6666
$ws.Range("K122").Value = 3094.1052
$ws.Range("L122").Value = 5665.9998
$ws.Range("M122").Value = -644.1052
$ws.Range("N122").Value = -10565.9998

$ws.Range("H126").Value = 1830.6666
$ws.Range("I126").Value = 1832.75
$ws.Range("J126").Value = 1814
$ws.Range("K126").Value = 5498.25
$ws.Range("L126").Value = 5442
$ws.Range("M126").Value = -3028.25
$ws.Range("N126").Value = -10382

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 2705.3784
$ws.Range("I100").Value = 1917.875
$ws.Range("J100").Value = 2922.6206
$ws.Range("K100").Value = 1917.875
$ws.Range("L100").Value = 2922.6206
$ws.Range("M100").Value = -1376.875
$ws.Range("N100").Value = -4004.6206

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 478267.66
$ws.Range("I122").Value = 834876.75
$ws.Range("J122").Value = 2788.889
$ws.Range("K122").Value = 2504630.25
$ws.Range("L122").Value = 8366.667000000001
$ws.Range("M122").Value = -2502180.25
$ws.Range("N122").Value = -13266.667
